$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 38

$ws.Cells.Item($row, 1).Value = "05/01/2026 03:39:38"
$ws.Cells.Item($row, 2).Value = "05/01 03:00"
$ws.Cells.Item($row, 3).Value = "g1 > Economia"
$ws.Cells.Item($row, 4).Value = "IPVA 2026: veja calendário de pagamento nos estados e no DF"
$ws.Cells.Item($row, 5).Value = "https://g1.globo.com/carros/noticia/2026/01/05/ipva-2026-veja-calendario-de-pagamento-nos-estados-e-no-df.ghtml"
$ws.Cells.Item($row, 6).Value = "congresso"
$ws.Cells.Item($row, 7).Value = "m a não pagar mais o IPVA, mas o período considerado varia.`nPorém, em dezembro de 2025, o Congresso Nacional promulgou uma proposta de emenda à Constituição (PEC) que isenta veículos fabric"
